# TP A3 du 13/02/2024
# Fill in the journal-de-bord row for the new session (row 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date of the session (stored as a real date serial, no time component).
$ws.Range("A14").Value = (Get-Date -Year 2024 -Month 2 -Day 13 -Hour 0 -Minute 0 -Second 0)

# Teacher / session type / group columns.
$ws.Range("B14").Value = "FSIL"
$ws.Range("C14").Value = "TP"
$ws.Range("F14").Value = "x"

# Description column.
$ws.Range("G14").Value = "QuestionsScore fix #2 : ce qui n'ont pas fini doivent finir pour la prochaine séance afin que tout le monde démarre la question 3"

# Commentaires column.
$ws.Range("I14").Value = "Groupe très hétérogène.
<!> Problème de licence serveur Intellij. Le SI est sur le coup mais cela ne fonctionne pas, il fut donc que les étudiants laisse le proxy activé, et choisisse un licence associée à un compte jetbrins qu'ils crééent si ils ne l'ont pas déjà avec l'adresse en ut-capitole pour avoir leur propre licence education..."

# Row grew taller once it held real (wrapped) text.
$ws.Rows.Item(14).RowHeight = 51.75
